# Update countries & provincias Spain
#
# This script applies the data refresh that was captured in the commit:
#  - Updates the "Datos actualizados..." timestamp footer (row 1).
#  - Refreshes COVID-19 stats for several countries whose ranking (by total
#    cases, column B, sorted descending) did not change (India, Banglades,
#    Emiratos Arabes Unidos, Suiza, Nepal, Madagascar, Senegal, Vietnam,
#    Liechtenstein).
#  - Refreshes Malta's stats; because Malta's new total (1470) now exceeds
#    Jordania's (1438), Malta moves above Jordania/Bahamas in the sorted
#    table, so rows 144-146 are rewritten as Malta, Jordania, Bahamas
#    (Jordania/Bahamas keep their previous figures, just shifted one row).
#  - Refreshes Montserrat/Islas Malvinas: Montserrat's active-cases count
#    overtakes Islas Malvinas, so the two swap places (rows 213-214).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 19 de Agosto de 2020 a las 13:27"

# --- In-place refreshes (country order unchanged) ----------------------

# Row 6: India
$ws.Cells.Item(6, 2).Value = 2771958
$ws.Cells.Item(6, 3).Value = 5332
$ws.Cells.Item(6, 4).Value = 2038709
$ws.Cells.Item(6, 5).Value = 680203
$ws.Cells.Item(6, 7).Value = 32
$ws.Cells.Item(6, 8).Value = 53046

# Row 19: Banglades
$ws.Cells.Item(19, 2).Value = 285091
$ws.Cells.Item(19, 3).Value = 2747
$ws.Cells.Item(19, 4).Value = 165738
$ws.Cells.Item(19, 5).Value = 115572
$ws.Cells.Item(19, 7).Value = 41
$ws.Cells.Item(19, 8).Value = 3781

# Row 44: Emiratos Arabes Unidos
$ws.Cells.Item(44, 2).Value = 65341
$ws.Cells.Item(44, 3).Value = 435
$ws.Cells.Item(44, 4).Value = 58022
$ws.Cells.Item(44, 5).Value = 6952
$ws.Cells.Item(44, 7).Value = 1
$ws.Cells.Item(44, 8).Value = 367

# Row 59: Suiza
$ws.Cells.Item(59, 2).Value = 38760
$ws.Cells.Item(59, 3).Value = 311
$ws.Cells.Item(59, 5).Value = 3264
$ws.Cells.Item(59, 7).Value = 4
$ws.Cells.Item(59, 8).Value = 1996

# Row 69: Nepal
$ws.Cells.Item(69, 2).Value = 28938
$ws.Cells.Item(69, 3).Value = 681
$ws.Cells.Item(69, 4).Value = 17700
$ws.Cells.Item(69, 5).Value = 11118
$ws.Cells.Item(69, 7).Value = 6
$ws.Cells.Item(69, 8).Value = 120

# Row 82: Madagascar
$ws.Cells.Item(82, 2).Value = 14074
$ws.Cells.Item(82, 3).Value = 65
$ws.Cells.Item(82, 4).Value = 12921
$ws.Cells.Item(82, 5).Value = 980

# Row 85: Senegal
$ws.Cells.Item(85, 2).Value = 12446
$ws.Cells.Item(85, 3).Value = 141
$ws.Cells.Item(85, 4).Value = 7877
$ws.Cells.Item(85, 5).Value = 4311
$ws.Cells.Item(85, 7).Value = 2
$ws.Cells.Item(85, 8).Value = 258

# Row 158: Vietnam
$ws.Cells.Item(158, 2).Value = 993
$ws.Cells.Item(158, 3).Value = 4
$ws.Cells.Item(158, 4).Value = 533
$ws.Cells.Item(158, 5).Value = 435
$ws.Cells.Item(158, 8).Value = 25

# Row 194: Liechtenstein
$ws.Cells.Item(194, 2).Value = 98
$ws.Cells.Item(194, 3).Value = 1
$ws.Cells.Item(194, 5).Value = 9

# --- Re-ranked block: Malta overtakes Jordania and Bahamas --------------
# Row 144 becomes Malta (new, higher figures)
$ws.Cells.Item(144, 1).Value = "Malta"
$ws.Cells.Item(144, 2).Value = 1470
$ws.Cells.Item(144, 3).Value = 47
$ws.Cells.Item(144, 4).Value = 784
$ws.Cells.Item(144, 5).Value = 677
$ws.Cells.Item(144, 8).Value = 9

# Row 145 becomes Jordania (figures unchanged, just shifted down a row)
$ws.Cells.Item(145, 1).Value = "Jordania"
$ws.Cells.Item(145, 2).Value = 1438
$ws.Cells.Item(145, 3).Value = 0
$ws.Cells.Item(145, 4).Value = 1243
$ws.Cells.Item(145, 5).Value = 184
$ws.Cells.Item(145, 8).Value = 11

# Row 146 becomes Bahamas (figures unchanged, just shifted down a row)
$ws.Cells.Item(146, 1).Value = "Bahamas"
$ws.Cells.Item(146, 2).Value = 1424
$ws.Cells.Item(146, 3).Value = 0
$ws.Cells.Item(146, 4).Value = 203
$ws.Cells.Item(146, 5).Value = 1201
$ws.Cells.Item(146, 8).Value = 20

# --- Re-ranked block: Montserrat overtakes Islas Malvinas ---------------
# Row 213 becomes Montserrat
$ws.Cells.Item(213, 1).Value = "Montserrat"
$ws.Cells.Item(213, 4).Value = 12
$ws.Cells.Item(213, 8).Value = 1

# Row 214 becomes Islas Malvinas
$ws.Cells.Item(214, 1).Value = "Islas Malvinas"
$ws.Cells.Item(214, 4).Value = 13
$ws.Cells.Item(214, 8).Value = 0
